$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for A2:C15 (sample, replicate, hours)
# Values below reflect the re-ordering / shifting of the time series so that the
# previous last row ("815a") becomes the new first data row, with every "hours"
# value re-based (offset) accordingly, and the old formula in C15 replaced by a
# plain numeric literal.
$data = @(
    @("815a",    "0", 0),
    @("1044a",   "0", 2.4833333333332988),
    @("1139a",   "0", 3.3999999999999657),
    @("1240p-1", "1", 4.4166666666666288),
    @("1240p",   "0", 4.4166666666666288),
    @("140p-1",  "1", 5.4166666666666288),
    @("140p",    "0", 5.4166666666666288),
    @("239p",    "0", 6.3999999999999684),
    @("338p-1",  "1", 7.3833333333332991),
    @("338p",    "0", 7.3833333333332991),
    @("454p",    "0", 8.6499999999999684),
    @("542p",    "0", 9.4499999999999691),
    @("640p",    "0", 10.416666666666629),
    @("743p",    "0", 11.466666666666629)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Clear the formula that used to live in the old C15 (now C15 holds a plain value)
$ws.Range("C15").Formula = $null
$ws.Range("C15").Value = 11.466666666666629

# Update the selected cell shown in the sheet view
$ws.Range("F6").Select()

$wb.Save()
